$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "Version: " + $newVersion
$wsAbout.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Leer South Coal Mine, United States, M1037, version ' + "'" + $newVersion + "'" + '. (See the CC license for attribution requirements if sharing or adapting the data set.)'

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 8; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
